# Atualização de bases das ligas, do dia: 29-02-2024 às 07:50
#
# The source data for rows 88-89, 95-96 and 104-107 on the
# "Estonia Meistriliiga" sheet had been associated with the wrong
# fixture (match id / teams / score / result / all odds columns travel
# together as a single record). This fixes the mismatch by
# re-assigning each record (columns B, F through AC) to its correct
# row while columns A (row index), C/D (league name) and E (match
# date) are left untouched since they were already correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88
$ws.Range("B88").Value = 6376945
$ws.Range("F88").Value = 'Parnu JK Vaprus'
$ws.Range("G88").Value = 'Harju JK Laagri'
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 'D'
$ws.Range("K88").Value = 1.615
$ws.Range("L88").Value = 4
$ws.Range("M88").Value = 4.5
$ws.Range("N88").Value = 1.85
$ws.Range("O88").Value = 3.8
$ws.Range("P88").Value = 3.5
$ws.Range("Q88").Value = -0.5
$ws.Range("R88").Value = 1.875
$ws.Range("S88").Value = 1.925
$ws.Range("T88").Value = 2.5
$ws.Range("U88").Value = 1.75
$ws.Range("V88").Value = 1.95
$ws.Range("W88").Value = -1
$ws.Range("X88").Value = 2.8
$ws.Range("Y88").Value = -1
$ws.Range("Z88").Value = -1
$ws.Range("AA88").Value = 0.925
$ws.Range("AB88").Value = -1
$ws.Range("AC88").Value = 0.95

# Row 89
$ws.Range("B89").Value = 6376947
$ws.Range("F89").Value = 'JK Tammeka Tartu'
$ws.Range("G89").Value = 'JK Tallinna Kalev'
$ws.Range("H89").Value = 2
$ws.Range("I89").Value = 7
$ws.Range("J89").Value = 'A'
$ws.Range("K89").Value = 3.6
$ws.Range("L89").Value = 3.4
$ws.Range("M89").Value = 1.909
$ws.Range("N89").Value = 2.4
$ws.Range("O89").Value = 3.6
$ws.Range("P89").Value = 2.45
$ws.Range("Q89").Value = 0
$ws.Range("R89").Value = 1.875
$ws.Range("S89").Value = 1.925
$ws.Range("T89").Value = 2.75
$ws.Range("U89").Value = 1.975
$ws.Range("V89").Value = 1.825
$ws.Range("W89").Value = -1
$ws.Range("X89").Value = -1
$ws.Range("Y89").Value = 1.45
$ws.Range("Z89").Value = -1
$ws.Range("AA89").Value = 0.925
$ws.Range("AB89").Value = 0.9750000000000001
$ws.Range("AC89").Value = -1

# Row 95
$ws.Range("B95").Value = 6416370
$ws.Range("F95").Value = 'FC Levadia Tallinn'
$ws.Range("G95").Value = 'Parnu JK Vaprus'
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 'D'
$ws.Range("K95").Value = 1.166
$ws.Range("L95").Value = 7
$ws.Range("M95").Value = 11
$ws.Range("N95").Value = 1.2
$ws.Range("O95").Value = 6
$ws.Range("P95").Value = 11
$ws.Range("Q95").Value = -2
$ws.Range("R95").Value = 1.85
$ws.Range("S95").Value = 1.95
$ws.Range("T95").Value = 3
$ws.Range("U95").Value = 1.85
$ws.Range("V95").Value = 1.95
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = 5
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = -1
$ws.Range("AA95").Value = 0.95
$ws.Range("AB95").Value = -1
$ws.Range("AC95").Value = 0.95

# Row 96
$ws.Range("B96").Value = 6482819
$ws.Range("F96").Value = 'JK Tammeka Tartu'
$ws.Range("G96").Value = 'FC Kuressaare'
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 'A'
$ws.Range("K96").Value = 1.833
$ws.Range("L96").Value = 3.5
$ws.Range("M96").Value = 3.5
$ws.Range("N96").Value = 2.1
$ws.Range("O96").Value = 3.4
$ws.Range("P96").Value = 2.875
$ws.Range("Q96").Value = -0.25
$ws.Range("R96").Value = 1.975
$ws.Range("S96").Value = 1.825
$ws.Range("T96").Value = 3
$ws.Range("U96").Value = 1.825
$ws.Range("V96").Value = 1.975
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = 1.875
$ws.Range("Z96").Value = -1
$ws.Range("AA96").Value = 0.825
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.9750000000000001

# Row 104
$ws.Range("B104").Value = 6533597
$ws.Range("F104").Value = 'FC Kuressaare'
$ws.Range("G104").Value = 'Parnu JK Vaprus'
$ws.Range("H104").Value = 1
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 'H'
$ws.Range("K104").Value = 2.5
$ws.Range("L104").Value = 3.4
$ws.Range("M104").Value = 2.5
$ws.Range("N104").Value = 2.15
$ws.Range("O104").Value = 3.6
$ws.Range("P104").Value = 2.875
$ws.Range("Q104").Value = -0.25
$ws.Range("R104").Value = 1.95
$ws.Range("S104").Value = 1.85
$ws.Range("T104").Value = 2.75
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = 1.85
$ws.Range("W104").Value = 1.15
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.95
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.8500000000000001

# Row 105
$ws.Range("B105").Value = 6537957
$ws.Range("F105").Value = 'FC Flora Tallinn'
$ws.Range("G105").Value = 'JK Nomme Kalju'
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 'D'
$ws.Range("K105").Value = 1.4
$ws.Range("L105").Value = 4
$ws.Range("M105").Value = 7.5
$ws.Range("N105").Value = 1.5
$ws.Range("O105").Value = 4.2
$ws.Range("P105").Value = 5
$ws.Range("Q105").Value = -1
$ws.Range("R105").Value = 1.85
$ws.Range("S105").Value = 1.95
$ws.Range("T105").Value = 2.75
$ws.Range("U105").Value = 1.85
$ws.Range("V105").Value = 1.95
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = 3.2
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 0.95
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0.95

# Row 106
$ws.Range("B106").Value = 6537869
$ws.Range("F106").Value = 'JK Tallinna Kalev'
$ws.Range("G106").Value = 'JK Trans Narva'
$ws.Range("H106").Value = 5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 'H'
$ws.Range("K106").Value = 1.6
$ws.Range("L106").Value = 4
$ws.Range("M106").Value = 4.5
$ws.Range("N106").Value = 1.65
$ws.Range("O106").Value = 4
$ws.Range("P106").Value = 4.333
$ws.Range("Q106").Value = -0.75
$ws.Range("R106").Value = 1.8
$ws.Range("S106").Value = 2
$ws.Range("T106").Value = 2.75
$ws.Range("U106").Value = 1.9
$ws.Range("V106").Value = 1.9
$ws.Range("W106").Value = 0.6499999999999999
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = -1
$ws.Range("Z106").Value = 0.8
$ws.Range("AA106").Value = -1
$ws.Range("AB106").Value = 0.8999999999999999
$ws.Range("AC106").Value = -1

# Row 107
$ws.Range("B107").Value = 6535416
$ws.Range("F107").Value = 'Paide Linnameeskond'
$ws.Range("G107").Value = 'FC Levadia Tallinn'
$ws.Range("H107").Value = 2
$ws.Range("I107").Value = 2
$ws.Range("J107").Value = 'D'
$ws.Range("K107").Value = 3
$ws.Range("L107").Value = 3.8
$ws.Range("M107").Value = 2
$ws.Range("N107").Value = 3
$ws.Range("O107").Value = 4
$ws.Range("P107").Value = 1.909
$ws.Range("Q107").Value = 0.5
$ws.Range("R107").Value = 1.85
$ws.Range("S107").Value = 1.95
$ws.Range("T107").Value = 2.75
$ws.Range("U107").Value = 1.95
$ws.Range("V107").Value = 1.85
$ws.Range("W107").Value = -1
$ws.Range("X107").Value = 3
$ws.Range("Y107").Value = -1
$ws.Range("Z107").Value = 0.8500000000000001
$ws.Range("AA107").Value = -1
$ws.Range("AB107").Value = 0.95
$ws.Range("AC107").Value = -1

